$wb = $excel.ActiveWorkbook

# Add the new "債務" (debt) sheet at the end of the workbook, after the
# existing last sheet ("保險").
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "債務"

# Header row (row 1) - bold/bordered header style carried over from the
# other sheets starts at column B, column A is left blank in the header.
$ws.Range("B1").Value = "species"
$ws.Range("C1").Value = "debtor"
$ws.Range("D1").Value = "owner"
$ws.Range("E1").Value = "total"
$ws.Range("F1").Value = "register_date"
$ws.Range("G1").Value = "register_reason"
$ws.Range("H1").Value = "property_category"
$ws.Range("I1").Value = "category"
$ws.Range("J1").Value = "date"
$ws.Range("K1").Value = "legislator_name"
$ws.Range("L1").Value = "legislator_id"
$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"

# Data row (row 2).
$ws.Range("A2").Value = 93
$ws.Range("B2").Value = "房屋貸款"
$ws.Range("C2").Value = "邱志偉"
$ws.Range("D2").Value = "中國信託商業銀行高雄分行"
$ws.Range("E2").Value = 1747472
$ws.Range("F2").Value = "94年07月27日"
$ws.Range("G2").Value = "購置房屋"
$ws.Range("H2").Value = "debt"
$ws.Range("I2").Value = "normal"
$ws.Range("J2").Value = "2012-03-06"
$ws.Range("K2").Value = "邱志偉"
$ws.Range("L2").Value = 1744
$ws.Range("M2").Value = "tmpba221"
$ws.Range("N2").Value = 93

# Style: bold/centered/bordered header row (matches the other sheets'
# "style 1"), plain default style for the data row ("style 2").
$ws.Range("A1:N1").Font.Bold = $true
$ws.Range("A1:N1").HorizontalAlignment = -4108
$ws.Range("A1:N1").VerticalAlignment = -4160
$ws.Range("A1:N1").Borders.LineStyle = 1

Write-Output "debt sheet added"
